# Applies the SLFRF Broadband Location Bulk Upload template update:
#  - insert a new "Location ID" column (C) used only when updating existing
#    bulk-upload entries
#  - append a new "Void Location" column (Q) used to void/display a location
#  - a few Required -> Conditional field changes
#  - updated help text for the Fabric ID column
#  - cosmetic view updates (zoom, column width, selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the new "Location ID" column before the existing Latitude column
#    (old column C), i.e. right after the Fabric ID column (B).
# ---------------------------------------------------------------------------
$ws.Columns("C:C").Insert()

$ws.Range("C4").Value = "Name"
$ws.Range("C5").Value = "Optional"
$ws.Range("C6").Value = "Name"
$ws.Range("C7").Value = "[DO NOT ENTER FOR NEW ENTRIES] This is a Treasury Portal auto-generated unique ID only for used when updating existing locations in bulk upload. The IDs can be found in the " + [char]0x201C + "My Broadband Locations" + [char]0x201D + " table as well as through the " + [char]0x201C + "Download as CSV" + [char]0x201D + " feature. "

$ws.Range("C4:C7").Font.Name = "Arial"
$ws.Range("C4:C7").Font.Size = 10
$ws.Range("C3").Font.Name = "Arial"
$ws.Range("C3").Font.Size = 10

# ---------------------------------------------------------------------------
# 2. Append the new "Void Location" column at the end (new column Q).
# ---------------------------------------------------------------------------
$ws.Range("Q4").Value = "Void_Location__c"
$ws.Range("Q5").Value = "Optional"
$ws.Range("Q6").Value = "Void Location"
$ws.Range("Q7").Value = "Not required. If the location is no longer relevant, set this column to 'Void' to remove it from the displayed locations. Locations will default to 'Display'." + [char]10 + "Options: " + [char]10 + "'Void'" + [char]10 + "'Display'"

# ---------------------------------------------------------------------------
# 3. Field-level changes: a few fields move from Required to Conditional, and
#    the Fabric ID help text is updated to flag it as an existing-entry-only
#    field.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Conditional"
$ws.Range("B7").Value = "[DO NOT ENTER FOR EXISTING ENTRIES]" + [char]10 + "Enter the site specific fabric identification number (Fabric ID) from the FCC broadband funded locations map that corresponds with the location in the field provided. This is the Broadband Serviceable Fabric Location (Max length - 20 characters)"

$ws.Range("M5").Value = "Conditional"
$ws.Range("O5").Value = "Conditional"
$ws.Range("P5").Value = "Conditional"

# ---------------------------------------------------------------------------
# 4. Cosmetic / view updates that came along with the re-save.
# ---------------------------------------------------------------------------
$ws.Columns("B:B").ColumnWidth = 35.78

$ws.Range("B8").Select()
$excel.ActiveWindow.Zoom = 80
